$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (refreshed stock-consistency data) ---
$ws.Range("C5").Value = -5
$ws.Range("D5").Value = 45849.29218469877
$ws.Range("E5").Value = -5
$ws.Range("F5").Value = 45848.7749537037
$ws.Range("C12").Value = 25
$ws.Range("D12").Value = 45849.29218469877
$ws.Range("E12").Value = 25
$ws.Range("F12").Value = 45848.68347222222
$ws.Range("C23").Value = 59
$ws.Range("D23").Value = 45849.29218469877
$ws.Range("E23").Value = 59
$ws.Range("F23").Value = 45848.49013888889
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 45849.29218469877
$ws.Range("E24").Value = 30
$ws.Range("F24").Value = 45848.62158564815
$ws.Range("C28").Value = 24
$ws.Range("D28").Value = 45849.29218469877
$ws.Range("E28").Value = 24
$ws.Range("F28").Value = 45848.73649305556
$ws.Range("C49").Value = 44
$ws.Range("D49").Value = 45849.29218469877
$ws.Range("E49").Value = 44
$ws.Range("F49").Value = 45848.7749537037
$ws.Range("C56").Value = 104
$ws.Range("D56").Value = 45849.29218469877
$ws.Range("E56").Value = 104
$ws.Range("F56").Value = 45848.69836805556
$ws.Range("C58").Value = 41
$ws.Range("D58").Value = 45849.29218469877
$ws.Range("E58").Value = 41
$ws.Range("F58").Value = 45848.64309027778
$ws.Range("C69").Value = 50
$ws.Range("D69").Value = 45849.29218469877
$ws.Range("E69").Value = 50
$ws.Range("F69").Value = 45848.7222800926
$ws.Range("C70").Value = 2
$ws.Range("D70").Value = 45849.29218469877
$ws.Range("E70").Value = 2
$ws.Range("F70").Value = 45848.73472222222
$ws.Range("C82").Value = 32
$ws.Range("D82").Value = 45849.29218469877
$ws.Range("E82").Value = 32
$ws.Range("F82").Value = 45848.64309027778
$ws.Range("C83").Value = 164
$ws.Range("D83").Value = 45849.29218469877
$ws.Range("E83").Value = 164
$ws.Range("F83").Value = 45848.7749537037
$ws.Range("C94").Value = 37
$ws.Range("D94").Value = 45849.29218469877
$ws.Range("E94").Value = 37
$ws.Range("F94").Value = 45848.7749537037
$ws.Range("C96").Value = 69
$ws.Range("D96").Value = 45849.29218469877
$ws.Range("E96").Value = 69
$ws.Range("F96").Value = 45848.48304398148
$ws.Range("C101").Value = 575
$ws.Range("D101").Value = 45849.29218469877
$ws.Range("E101").Value = 575
$ws.Range("F101").Value = 45848.7749537037
$ws.Range("C125").Value = 190
$ws.Range("D125").Value = 45849.29218469877
$ws.Range("E125").Value = 190
$ws.Range("F125").Value = 45848.70797453704
$ws.Range("C139").Value = 70
$ws.Range("D139").Value = 45849.29218469877
$ws.Range("E139").Value = 70
$ws.Range("F139").Value = 45848.48304398148
$ws.Range("C141").Value = 223
$ws.Range("D141").Value = 45849.29218469877
$ws.Range("E141").Value = 223
$ws.Range("F141").Value = 45848.70797453704
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 45849.29218469877
$ws.Range("E146").Value = 0
$ws.Range("F146").Value = 45848.74967592592
$ws.Range("F156").Value = 45848.72150462963
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 45849.29218469877
$ws.Range("E181").Value = 0
$ws.Range("F181").Value = 45848.74856481481
$ws.Range("C195").Value = -7
$ws.Range("D195").Value = 45849.29218469877
$ws.Range("E195").Value = -7
$ws.Range("F195").Value = 45848.68099537037
$ws.Range("C200").Value = 938
$ws.Range("D200").Value = 45849.29218469877
$ws.Range("E200").Value = 938
$ws.Range("F200").Value = 45848.70797453704
$ws.Range("C201").Value = 106
$ws.Range("D201").Value = 45849.29218469877
$ws.Range("E201").Value = 106
$ws.Range("F201").Value = 45848.70797453704
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 45849.2922031681
$ws.Range("E206").Value = 0
$ws.Range("F206").Value = 45848.77842592593
$ws.Range("C217").Value = 18
$ws.Range("D217").Value = 45849.29218469877
$ws.Range("E217").Value = 18
$ws.Range("F217").Value = 45848.64112268519
$ws.Range("C270").Value = -2
$ws.Range("D270").Value = 45849.29218469877
$ws.Range("E270").Value = -2
$ws.Range("F270").Value = 45848.75108796296
$ws.Range("C281").Value = 24
$ws.Range("D281").Value = 45849.29218469877
$ws.Range("E281").Value = 24
$ws.Range("F281").Value = 45848.7749537037
$ws.Range("C295").Value = 26
$ws.Range("D295").Value = 45849.29218469877
$ws.Range("E295").Value = 26
$ws.Range("F295").Value = 45848.70797453704
$ws.Range("C300").Value = 10
$ws.Range("D300").Value = 45849.29218469877
$ws.Range("E300").Value = 10
$ws.Range("F300").Value = 45848.62158564815
$ws.Range("C309").Value = 884
$ws.Range("D309").Value = 45849.29218469877
$ws.Range("E309").Value = 884
$ws.Range("F309").Value = 45848.68347222222
$ws.Range("C346").Value = 47
$ws.Range("D346").Value = 45849.29218469877
$ws.Range("E346").Value = 47
$ws.Range("F346").Value = 45848.7749537037
$ws.Range("C351").Value = 990
$ws.Range("D351").Value = 45849.29218469877
$ws.Range("E351").Value = 990
$ws.Range("F351").Value = 45848.48304398148
$ws.Range("C363").Value = 297
$ws.Range("D363").Value = 45849.29218469877
$ws.Range("E363").Value = 297
$ws.Range("F363").Value = 45848.7749537037
$ws.Range("F377").Value = 45848.72100694444
$ws.Range("C387").Value = 9
$ws.Range("D387").Value = 45849.29218469877
$ws.Range("E387").Value = 9
$ws.Range("F387").Value = 45848.62158564815
$ws.Range("C390").Value = 115
$ws.Range("D390").Value = 45849.29218469877
$ws.Range("E390").Value = 115
$ws.Range("F390").Value = 45848.62158564815
$ws.Range("C394").Value = 96
$ws.Range("D394").Value = 45849.29218469877
$ws.Range("E394").Value = 96
$ws.Range("F394").Value = 45848.44458333333
$ws.Range("C404").Value = 402
$ws.Range("D404").Value = 45849.29218469877
$ws.Range("E404").Value = 402
$ws.Range("F404").Value = 45848.7749537037
$ws.Range("C422").Value = 4
$ws.Range("D422").Value = 45849.29218469877
$ws.Range("E422").Value = 4
$ws.Range("F422").Value = 45848.70797453704
$ws.Range("C441").Value = 2
$ws.Range("D441").Value = 45849.29218469877
$ws.Range("E441").Value = 2
$ws.Range("F441").Value = 45848.68099537037
$ws.Range("C442").Value = -5
$ws.Range("D442").Value = 45849.29218469877
$ws.Range("E442").Value = -5
$ws.Range("F442").Value = 45848.68099537037
$ws.Range("C472").Value = 8
$ws.Range("D472").Value = 45849.29218469877
$ws.Range("E472").Value = 8
$ws.Range("F472").Value = 45848.48304398148
$ws.Range("D488").Value = 45849.29218469877
$ws.Range("D489").Value = 45849.29218469877
$ws.Range("D490").Value = 45849.29218469877
$ws.Range("C510").Value = 109
$ws.Range("D510").Value = 45849.29218469877
$ws.Range("E510").Value = 109
$ws.Range("F510").Value = 45848.64309027778
$ws.Range("C570").Value = 2429
$ws.Range("D570").Value = 45849.29218469877
$ws.Range("E570").Value = 2429
$ws.Range("F570").Value = 45848.62158564815
$ws.Range("C582").Value = 17
$ws.Range("D582").Value = 45849.2922031681
$ws.Range("E582").Value = 17
$ws.Range("F582").Value = 45848.7749537037
$ws.Range("C596").Value = 4
$ws.Range("D596").Value = 45849.29218469877
$ws.Range("E596").Value = 4
$ws.Range("F596").Value = 45848.48304398148
$ws.Range("C631").Value = 21
$ws.Range("D631").Value = 45849.29218469877
$ws.Range("E631").Value = 21
$ws.Range("F631").Value = 45848.7075462963
$ws.Range("C657").Value = 1637
$ws.Range("D657").Value = 45849.2922031681
$ws.Range("E657").Value = 1637
$ws.Range("F657").Value = 45848.7749537037
$ws.Range("C680").Value = 0
$ws.Range("D680").Value = 45849.29218469877
$ws.Range("E680").Value = 0
$ws.Range("F680").Value = 45848.48304398148
$ws.Range("C716").Value = 17
$ws.Range("D716").Value = 45849.29218469877
$ws.Range("E716").Value = 17
$ws.Range("F716").Value = 45848.62158564815
$ws.Range("C726").Value = 19
$ws.Range("D726").Value = 45849.29218469877
$ws.Range("E726").Value = 19
$ws.Range("F726").Value = 45848.66717592593
$ws.Range("C729").Value = -4
$ws.Range("D729").Value = 45849.29218469877
$ws.Range("E729").Value = -4
$ws.Range("F729").Value = 45848.6952662037
$ws.Range("C821").Value = 84
$ws.Range("D821").Value = 45849.2922031681
$ws.Range("E821").Value = 84
$ws.Range("F821").Value = 45848.7749537037
$ws.Range("C826").Value = 152
$ws.Range("D826").Value = 45849.2922031681
$ws.Range("E826").Value = 152
$ws.Range("F826").Value = 45848.7749537037
$ws.Range("C844").Value = 1
$ws.Range("D844").Value = 45849.29218469877
$ws.Range("E844").Value = 1
$ws.Range("F844").Value = 45848.75552083334
$ws.Range("F845").Value = 45848.46533564815
$ws.Range("C883").Value = 248
$ws.Range("D883").Value = 45849.2922031681
$ws.Range("E883").Value = 248
$ws.Range("F883").Value = 45848.7749537037
$ws.Range("C900").Value = 2
$ws.Range("D900").Value = 45849.29218469877
$ws.Range("E900").Value = 2
$ws.Range("F900").Value = 45848.44458333333
$ws.Range("C1015").Value = -1
$ws.Range("D1015").Value = 45849.29218469877
$ws.Range("E1015").Value = -1
$ws.Range("F1015").Value = 45848.40421296296
$ws.Range("C1027").Value = 20
$ws.Range("D1027").Value = 45849.29218469877
$ws.Range("E1027").Value = 20
$ws.Range("F1027").Value = 45848.64112268519
$ws.Range("C1043").Value = 20
$ws.Range("D1043").Value = 45849.2922031681
$ws.Range("E1043").Value = 20
$ws.Range("F1043").Value = 45848.7749537037
$ws.Range("C1057").Value = 5
$ws.Range("D1057").Value = 45849.29218469877
$ws.Range("E1057").Value = 5
$ws.Range("F1057").Value = 45848.48304398148
$ws.Range("C1135").Value = 2
$ws.Range("D1135").Value = 45849.29218469877
$ws.Range("E1135").Value = 2
$ws.Range("F1135").Value = 45848.46533564815
$ws.Range("C1150").Value = 34
$ws.Range("D1150").Value = 45849.29218469877
$ws.Range("E1150").Value = 34
$ws.Range("F1150").Value = 45848.62158564815
$ws.Range("D1171").Value = 45849.29218469877
$ws.Range("C1193").Value = 4
$ws.Range("D1193").Value = 45849.29218469877
$ws.Range("E1193").Value = 4
$ws.Range("F1193").Value = 45848.70797453704
$ws.Range("C1253").Value = 799
$ws.Range("D1253").Value = 45849.2922031681
$ws.Range("E1253").Value = 799
$ws.Range("F1253").Value = 45848.7749537037
$ws.Range("C1332").Value = -5
$ws.Range("D1332").Value = 45849.29218469877
$ws.Range("E1332").Value = -5
$ws.Range("F1332").Value = 45848.73472222222
$ws.Range("C1342").Value = 831
$ws.Range("D1342").Value = 45849.29218469877
$ws.Range("E1342").Value = 831
$ws.Range("F1342").Value = 45848.70797453704
$ws.Range("C1345").Value = -4
$ws.Range("D1345").Value = 45849.29218469877
$ws.Range("E1345").Value = -4
$ws.Range("F1345").Value = 45848.7075462963
$ws.Range("C1393").Value = 26
$ws.Range("D1393").Value = 45849.29218469877
$ws.Range("E1393").Value = 26
$ws.Range("F1393").Value = 45848.66505787037
$ws.Range("C1403").Value = 36
$ws.Range("D1403").Value = 45849.29218469877
$ws.Range("E1403").Value = 36
$ws.Range("F1403").Value = 45848.70797453704
$ws.Range("C1411").Value = 58
$ws.Range("D1411").Value = 45849.29218469877
$ws.Range("E1411").Value = 58
$ws.Range("F1411").Value = 45848.70797453704
$ws.Range("C1423").Value = -10
$ws.Range("D1423").Value = 45849.29218469877
$ws.Range("E1423").Value = -10
$ws.Range("F1423").Value = 45848.7075462963
$ws.Range("C1427").Value = 0
$ws.Range("D1427").Value = 45849.29218469877
$ws.Range("E1427").Value = 0
$ws.Range("F1427").Value = 45848.7075462963
$ws.Range("C1448").Value = 15
$ws.Range("D1448").Value = 45849.29218469877
$ws.Range("E1448").Value = 15
$ws.Range("F1448").Value = 45848.70797453704
$ws.Range("C1507").Value = 140
$ws.Range("D1507").Value = 45849.2922031681
$ws.Range("E1507").Value = 140
$ws.Range("F1507").Value = 45848.7749537037
$ws.Range("C1512").Value = 10
$ws.Range("D1512").Value = 45849.29218469877
$ws.Range("E1512").Value = 10
$ws.Range("F1512").Value = 45848.74804398148
$ws.Range("C1567").Value = 0
$ws.Range("D1567").Value = 45849.29218469877
$ws.Range("E1567").Value = 0
$ws.Range("F1567").Value = 45848.48969907407
$ws.Range("C1594").Value = 120
$ws.Range("D1594").Value = 45849.29218469877
$ws.Range("E1594").Value = 120
$ws.Range("F1594").Value = 45848.62158564815
$ws.Range("C1597").Value = 3878
$ws.Range("D1597").Value = 45849.29218469877
$ws.Range("E1597").Value = 3878
$ws.Range("F1597").Value = 45848.70797453704
$ws.Range("C1844").Value = 345
$ws.Range("D1844").Value = 45849.2922031681
$ws.Range("E1844").Value = 345
$ws.Range("F1844").Value = 45848.7749537037
$ws.Range("C1880").Value = 0
$ws.Range("D1880").Value = 45849.29218469877
$ws.Range("E1880").Value = 0
$ws.Range("F1880").Value = 45848.64112268519
$ws.Range("C1957").Value = 57
$ws.Range("D1957").Value = 45849.2922031681
$ws.Range("E1957").Value = 57
$ws.Range("F1957").Value = 45848.7749537037
$ws.Range("C2023").Value = 61
$ws.Range("D2023").Value = 45849.2922031681
$ws.Range("E2023").Value = 61
$ws.Range("F2023").Value = 45848.7749537037
$ws.Range("C2024").Value = 55
$ws.Range("D2024").Value = 45849.29218469877
$ws.Range("E2024").Value = 55
$ws.Range("F2024").Value = 45848.62158564815
$ws.Range("C2299").Value = 2
$ws.Range("D2299").Value = 45849.29218469877
$ws.Range("E2299").Value = 2
$ws.Range("F2299").Value = 45848.62158564815
$ws.Range("C2375").Value = 0
$ws.Range("D2375").Value = 45849.29218469877
$ws.Range("E2375").Value = 0
$ws.Range("F2375").Value = 45848.71688657408

# --- Append new rows for newly tracked products ---
$ws.Range("A2609").Value = 44118271
$ws.Range("B2609").Value = 1
$ws.Range("C2609").Value = 0
$ws.Range("D2609").Value = 45849.29218469877
$ws.Range("D2609").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2609").Value = 0
$ws.Range("G2609").Value = 0
$ws.Range("H2609").Value = "Consistente"

$ws.Range("A2610").Value = 44119087
$ws.Range("B2610").Value = 1
$ws.Range("C2610").Value = 0
$ws.Range("D2610").Value = 45849.29218469877
$ws.Range("D2610").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2610").Value = 0
$ws.Range("G2610").Value = 0
$ws.Range("H2610").Value = "Consistente"

$ws.Range("A2611").Value = 44121363
$ws.Range("B2611").Value = 1
$ws.Range("C2611").Value = 0
$ws.Range("D2611").Value = 45849.29218469877
$ws.Range("D2611").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2611").Value = 0
$ws.Range("G2611").Value = 0
$ws.Range("H2611").Value = "Consistente"

$ws.Range("A2612").Value = 44121568
$ws.Range("B2612").Value = 1
$ws.Range("C2612").Value = 0
$ws.Range("D2612").Value = 45849.29218469877
$ws.Range("D2612").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2612").Value = 0
$ws.Range("G2612").Value = 0
$ws.Range("H2612").Value = "Consistente"

$ws.Range("A2613").Value = 44122289
$ws.Range("B2613").Value = 1
$ws.Range("C2613").Value = 0
$ws.Range("D2613").Value = 45849.29218469877
$ws.Range("D2613").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2613").Value = 0
$ws.Range("G2613").Value = 0
$ws.Range("H2613").Value = "Consistente"

